$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AUK0451", "ARIOVALDO SOUZA GOMES", 620578),
    @("DDY4C74", "MARIO ROBERTO", 620686),
    @("DPE0B20", "DAVID DE JESUS", 619353),
    @("EAR7C31", "LEONARDO MAGALHAES", 620867),
    @("GHG7C42", "RODOLFO PIZANI", 436321),
    @("IWJ4B20", "DAVID DE JESUS", 619353)
)

$startRow = 46
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
